$wb = $excel.ActiveWorkbook

# --- "Data" sheet: append the newest weekly WALCL observation ---
$dataSheet = $wb.Worksheets.Item("Data")

$lastRow = 111
$newRow = 112

# Copy the formatting (date number format / font / border / alignment) from the
# preceding date cell so the new row matches the existing look, then fill values.
$dataSheet.Range("A$lastRow").Copy()
$dataSheet.Range("A$newRow").PasteSpecial(-4122)

$dataSheet.Cells.Item($newRow, 1).Value = 45245
$dataSheet.Cells.Item($newRow, 2).Value = 7814.991

# --- "SeriesInfo" sheet: refresh metadata for the new data pull ---
$infoSheet = $wb.Worksheets.Item("SeriesInfo")

function Set-TextValue($range, $text) {
    # Force the assignment to be kept as literal text instead of letting Excel
    # auto-parse date-looking strings into date serial numbers, then drop the
    # temporary number format so the cell keeps its original (default) style.
    $range.NumberFormat = "@"
    $range.Value = $text
    $range.ClearFormats()
}

Set-TextValue $infoSheet.Range("B3") "2023-11-21"
Set-TextValue $infoSheet.Range("B4") "2023-11-21"
Set-TextValue $infoSheet.Range("B7") "2023-11-15"
Set-TextValue $infoSheet.Range("B14") "2023-11-16 15:33:02-06"
